$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must stay literal text even though it looks
# like a number/percentage. A leading apostrophe forces Excel to store it
# as text (quotePrefix) without altering the cell's number format.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

Set-TextValue "D2" "331.76"

Set-TextValue "E3" "1.66%"

Set-TextValue "D4" "5.683"
Set-TextValue "E4" "-4.90%"

Set-TextValue "D5" "0.08084"
Set-TextValue "E5" "-0.70%"

Set-TextValue "D6" "2.040"
Set-TextValue "E6" "3.62%"

Set-TextValue "D7" "8.739"
Set-TextValue "E7" "-0.29%"

Set-TextValue "D8" "4.534"
Set-TextValue "E8" "-1.71%"

Set-TextValue "E9" "0.93%"

Set-TextValue "D10" "0.9236"
Set-TextValue "E10" "-2.65%"

Set-TextValue "D11" "0.1261"
Set-TextValue "E11" "-4.55%"

Set-TextValue "D12" "0.1945"
Set-TextValue "E12" "-2.43%"

Set-TextValue "D13" "8.733"
Set-TextValue "E13" "-2.81%"

Set-TextValue "D14" "0.09333"
Set-TextValue "E14" "-0.43%"

Set-TextValue "D15" "0.03742"
Set-TextValue "E15" "6.98%"

Set-TextValue "E16" "9.38%"

Set-TextValue "D17" "0.001298"
Set-TextValue "E17" "-1.57%"

Set-TextValue "D18" "0.006302"
Set-TextValue "E18" "-1.33%"

Set-TextValue "D19" "3.365"
Set-TextValue "E19" "0.26%"

Set-TextValue "E20" "-1.76%"

Set-TextValue "D21" "0.1419"
Set-TextValue "E21" "0.37%"

Set-TextValue "D22" "0.2656"
Set-TextValue "E22" "9.97%"

Set-TextValue "D23" "0.04445"
Set-TextValue "E23" "0.07%"

Set-TextValue "D24" "0.001262"
Set-TextValue "E24" "-0.17%"

Set-TextValue "D25" "0.004293"
Set-TextValue "E25" "-3.31%"

Set-TextValue "D26" "0.0001242"
Set-TextValue "E26" "13.49%"

Set-TextValue "D39" "0.02859"
Set-TextValue "E39" "15.62%"

Set-TextValue "D40" "0.05474"
Set-TextValue "E40" "3.36%"

Set-TextValue "D41" "0.007794"
Set-TextValue "E41" "3.59%"

Set-TextValue "D42" "0.009945"
Set-TextValue "E42" "10.80%"

Set-TextValue "D43" "0.1420"
Set-TextValue "E43" "-0.99%"

Set-TextValue "D44" "0.002134"
Set-TextValue "E44" "3.63%"

Set-TextValue "D45" "0.01179"
Set-TextValue "E45" "14.23%"

Set-TextValue "D46" "0.00006778"

Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "-0.23%"

Set-TextValue "D48" "0.002283"
Set-TextValue "E48" "26.33%"

Set-TextValue "D49" "0.003023"
Set-TextValue "E49" "-13.85%"

Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "-0.23%"

Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "-0.23%"
